$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 52,5
$data[0,0] = 39400; $data[0,1] = 2007; $data[0,2] = 7.226520411029069; $data[0,3] = 2008; $data[0,4] = $null
$data[1,0] = 39583; $data[1,1] = 2008; $data[1,2] = $null; $data[1,3] = 2009; $data[1,4] = $null
$data[2,0] = 39765; $data[2,1] = 2008; $data[2,2] = 4.268860212333636; $data[2,3] = 2009; $data[2,4] = $null
$data[3,0] = 39948; $data[3,1] = 2009; $data[3,2] = $null; $data[3,3] = 2010; $data[3,4] = $null
$data[4,0] = 40130; $data[4,1] = 2009; $data[4,2] = -7.266312015249776; $data[4,3] = 2010; $data[4,4] = $null
$data[5,0] = 40310; $data[5,1] = 2010; $data[5,2] = 3.184002331674129; $data[5,3] = 2011; $data[5,4] = 6.704254199558113
$data[6,0] = 40494; $data[6,1] = 2010; $data[6,2] = 6.958243460951929; $data[6,3] = 2011; $data[6,4] = 12.21658306395068
$data[7,0] = 40676; $data[7,1] = 2011; $data[7,2] = 8.626810748872327; $data[7,3] = 2012; $data[7,4] = 4.739201070534826
$data[8,0] = 40862; $data[8,1] = 2011; $data[8,2] = 9.469137444079934; $data[8,3] = 2012; $data[8,4] = 8.079264579851909
$data[9,0] = 41044; $data[9,1] = 2012; $data[9,2] = 3.449685446853534; $data[9,3] = 2013; $data[9,4] = 3.26507595662513
$data[10,0] = 41228; $data[10,1] = 2012; $data[10,2] = 3.358206407534947; $data[10,3] = 2013; $data[10,4] = 4.701432377325987
$data[11,0] = 41409; $data[11,1] = 2013; $data[11,2] = -1.480934717826909; $data[11,3] = 2014; $data[11,4] = 0.7772706050320544
$data[12,0] = 41592; $data[12,1] = 2013; $data[12,2] = 0.3081076735359067; $data[12,3] = 2014; $data[12,4] = 3.972902167062387
$data[13,0] = 41774; $data[13,1] = 2014; $data[13,2] = 5.427992542801308; $data[13,3] = 2015; $data[13,4] = 4.945882057432871
$data[14,0] = 41957; $data[14,1] = 2014; $data[14,2] = 3.901355411819707; $data[14,3] = 2015; $data[14,4] = 4.658857392675264
$data[15,0] = 42137; $data[15,1] = 2015; $data[15,2] = 4.970284184513551; $data[15,3] = 2016; $data[15,4] = 4.488174889976171
$data[16,0] = 42321; $data[16,1] = 2015; $data[16,2] = 5.331683351557981; $data[16,3] = 2016; $data[16,4] = 4.089819750351786
$data[17,0] = 42503; $data[17,1] = 2016; $data[17,2] = 4.039484738713828; $data[17,3] = 2017; $data[17,4] = 4.214976960249173
$data[18,0] = 42689; $data[18,1] = 2016; $data[18,2] = 3.254758369308375; $data[18,3] = 2017; $data[18,4] = 2.313009565865753
$data[19,0] = 42867; $data[19,1] = 2017; $data[19,2] = 4.589070866863865; $data[19,3] = 2018; $data[19,4] = 3.829046580278361
$data[20,0] = 43053; $data[20,1] = 2017; $data[20,2] = 5.246209615995667; $data[20,3] = 2018; $data[20,4] = 4.784022165496182
$data[21,0] = 43145; $data[21,1] = 2018; $data[21,2] = 6.011890504679696; $data[21,3] = 2019; $data[21,4] = 4.234360353587641
$data[22,0] = 43235; $data[22,1] = 2018; $data[22,2] = 3.625873842174787; $data[22,3] = 2019; $data[22,4] = 2.330842103296149
$data[23,0] = 43326; $data[23,1] = 2018; $data[23,2] = 4.899902276557011; $data[23,3] = 2019; $data[23,4] = 3.75342745845737
$data[24,0] = 43418; $data[24,1] = 2018; $data[24,2] = 4.86255966374296; $data[24,3] = 2019; $data[24,4] = 4.112897401876747
$data[25,0] = 43510; $data[25,1] = 2019; $data[25,2] = 3.660106318836931; $data[25,3] = 2020; $data[25,4] = 3.270208315717005
$data[26,0] = 43600; $data[26,1] = 2019; $data[26,2] = 3.500574054404404; $data[26,3] = 2020; $data[26,4] = 3.21661481720994
$data[27,0] = 43691; $data[27,1] = 2019; $data[27,2] = 2.983312281417039; $data[27,3] = 2020; $data[27,4] = 2.428295356218069
$data[28,0] = 43783; $data[28,1] = 2019; $data[28,2] = 2.764740011159428; $data[28,3] = 2020; $data[28,4] = 1.643374185611401
$data[29,0] = 43875; $data[29,1] = 2020; $data[29,2] = 2.096953540210977; $data[29,3] = 2021; $data[29,4] = 3.169670668618951
$data[30,0] = 43966; $data[30,1] = 2020; $data[30,2] = -0.9913189363815245; $data[30,3] = 2021; $data[30,4] = 1.183532150252908
$data[31,0] = 44068; $data[31,1] = 2020; $data[31,2] = -7.578477024949737; $data[31,3] = 2021; $data[31,4] = -5.743787238149123
$data[32,0] = 44159; $data[32,1] = 2020; $data[32,2] = -7.260793671746435; $data[32,3] = 2021; $data[32,4] = 0.00562230452727519
$data[33,0] = 44251; $data[33,1] = 2021; $data[33,2] = 0.4989366167094333; $data[33,3] = 2022; $data[33,4] = 2.690694906265412
$data[34,0] = 44341; $data[34,1] = 2021; $data[34,2] = 4.507091823899212; $data[34,3] = 2022; $data[34,4] = 5.429743376942153
$data[35,0] = 44432; $data[35,1] = 2021; $data[35,2] = 4.379227219808146; $data[35,3] = 2022; $data[35,4] = 4.954652839642848
$data[36,0] = 44525; $data[36,1] = 2021; $data[36,2] = 4.097586525396268; $data[36,3] = 2022; $data[36,4] = 3.9116372951149
$data[37,0] = 44617; $data[37,1] = 2022; $data[37,2] = 7.041577295022128; $data[37,3] = 2023; $data[37,4] = 3.388682041315016
$data[38,0] = 44706; $data[38,1] = 2022; $data[38,2] = 8.053468068361846; $data[38,3] = 2023; $data[38,4] = 3.974997080343634
$data[39,0] = 44798; $data[39,1] = 2022; $data[39,2] = 7.397318165265498; $data[39,3] = 2023; $data[39,4] = 3.367096865515662
$data[40,0] = 44890; $data[40,1] = 2022; $data[40,2] = 7.824284864703746; $data[40,3] = 2023; $data[40,4] = 2.586378346096296
$data[41,0] = 44981; $data[41,1] = 2023; $data[41,2] = 0.9995490351194292; $data[41,3] = 2024; $data[41,4] = 2.834404338648921
$data[42,0] = 45071; $data[42,1] = 2023; $data[42,2] = 0.2714278794373248; $data[42,3] = 2024; $data[42,4] = 2.460471645027118
$data[43,0] = 45163; $data[43,1] = 2023; $data[43,2] = -0.3046246622258053; $data[43,3] = 2024; $data[43,4] = 1.976476469605681
$data[44,0] = 45254; $data[44,1] = 2023; $data[44,2] = -1.24502235313334; $data[44,3] = 2024; $data[44,4] = -1.561801765212567
$data[45,0] = 45345; $data[45,1] = 2024; $data[45,2] = -2.798317913999848; $data[45,3] = 2025; $data[45,4] = 2.363509743917169
$data[46,0] = 45436; $data[46,1] = 2024; $data[46,2] = -2.107534670984712; $data[46,3] = 2025; $data[46,4] = 2.747596279389564
$data[47,0] = 45534; $data[47,1] = 2024; $data[47,2] = -2.567041707495976; $data[47,3] = 2025; $data[47,4] = 1.835066812373642
$data[48,0] = 45618; $data[48,1] = 2024; $data[48,2] = -1.735114423676209; $data[48,3] = 2025; $data[48,4] = 2.409056355286521
$data[49,0] = 45713; $data[49,1] = 2025; $data[49,2] = 2.450219408996213; $data[49,3] = 2026; $data[49,4] = 2.677741483899121
$data[50,0] = 45800; $data[50,1] = 2025; $data[50,2] = 1.552685227480533; $data[50,3] = 2026; $data[50,4] = 2.496145622272206
$data[51,0] = 45891; $data[51,1] = 2025; $data[51,2] = 2.64031107104763; $data[51,3] = 2026; $data[51,4] = 3.123685491361705

# Write the full updated dataset (header row untouched) in one shot
$ws.Range("A2:E53").Value = $data

# Row 53 is a brand-new row beyond the old A1:E52 used range, so it has no
# inherited style; copy the date-format style from the row above (A52) onto
# A53 so column A keeps its "YYYY-MM-DD HH:MM:SS" display like every other row.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

